$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Sr. No" header in A1 to "Face_ID"
$ws.Range("A1").Value = "Face_ID"

# Update the active selection to A2 (was E9)
$ws.Range("A2").Select()
